$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$sh = $s.Shapes.Item("TextBox 17")

# Shrink the textbox: width 2540000 EMU -> 2270328 EMU (height stays 1169551 EMU)
$sh.Width = 2270328 / 12700

# First run of the text ("Time normalisation ") switches from bold-italic
# "CMU Serif BoldItalic" to bold-only "CMU Serif"
$tr = $sh.TextFrame.TextRange
$run1 = $tr.Characters(1, 19)
$run1.Font.Italic = $false
$run1.Font.Name = "CMU Serif"
$run1.Font.NameFarEast = "CMU Serif"
$run1.Font.NameComplexScript = "CMU Serif"
